$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leakage Bug List")
$cell = $ws.Range("A1")
$cell.Value = "Hello World"
$cell.Characters(1, 5).Font.Bold = $true
$len = $cell.Characters().Count
$cell.Characters($len+1, 0).Text = "!!!"
# reapply bold to first 5 chars since it got flattened
$cell.Characters(1, 5).Font.Bold = $true
Write-Host "after re-apply: " $cell.Value2
